# DOMA-1030 - refactor with global mappers
#
# Applies the TicketAnalyticsExportTemplate[status_property].xlsx changes:
#  1. Rename the "Complete" status label to "Completed" (cell C1).
#  2. Move the saved cell selection from D9 to E9.
#  3. Widen/resize columns B:E to their new widths.
#  4. Shrink the header row heights (rows 2 and 3) from 36.9 to 25.1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. "Complete" -> "Completed"
$ws.Range("C1").Value = "Completed"

# 2. Selection: D9 -> E9
$ws.Range("E9").Select() | Out-Null

# 3. Column widths (character units, same scale as the legacy ColumnWidth property)
$ws.Columns.Item(2).ColumnWidth = 23.166666666666668
$ws.Columns.Item(3).ColumnWidth = 16.166666666666668
$ws.Columns.Item(4).ColumnWidth = 15.333333333333334
$ws.Columns.Item(5).ColumnWidth = 13.833333333333334

# 4. Row heights for the two data rows
$ws.Rows.Item(2).RowHeight = 25.1
$ws.Rows.Item(3).RowHeight = 25.1

Write-Host "Applied DOMA-1030 template tweaks"
